$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G: rename header from "image_name" to "views", and replace every
# row's "a.png" text value with a descending numeric view count.
$ws.Range("G1").Value = "views"

$views = @{
    2  = 14
    3  = 13
    4  = 12
    5  = 11
    6  = 10
    7  = 9
    8  = 8
    9  = 7
    10 = 6
    11 = 5
    12 = 4
    13 = 3
    14 = 2
    15 = 1
}
foreach ($row in $views.Keys) {
    $ws.Cells.Item($row, 7).Value = $views[$row]
}

# Column D: shorten street names.
$ws.Range("D4").Value = "Gr Western Rd, Glasgow"
$ws.Range("D5").Value = "Gr Western Rd, Glasgow"
$ws.Range("D6").Value = "Gr Western Rd, Glasgow"

$ws.Range("D8").Value = "Buchanan St, Glasgow"
$ws.Range("D9").Value = "Buchanan St, Glasgow"

$ws.Range("D10").Value = "Shuna St, Glasgow"
$ws.Range("D11").Value = "Shuna St, Glasgow"
$ws.Range("D12").Value = "Shuna St, Glasgow"

# Update the selected cell in the sheet view.
$ws.Range("D2").Select()
